$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("U2").Value = 1.93
$ws.Range("Y2").Value = 15
$ws.Range("AB2").Value = 8
$ws.Range("F3").Value = 1.02
$ws.Range("G3").Value = 980
$ws.Range("H3").Value = 1.01
$ws.Range("I3").Value = 980
$ws.Range("J3").Value = 1.2
$ws.Range("K3").Value = 980
$ws.Range("N3").Value = 1.24
$ws.Range("P3").Value = 1.24
$ws.Range("Q3").Value = 1.44
$ws.Range("S3").Value = 1.44
$ws.Range("T3").Value = 1.04
$ws.Range("U3").Value = 1.04
$ws.Range("V3").Value = 1.01
$ws.Range("W3").Value = 1.8
$ws.Range("F4").Value = 1.02
$ws.Range("G4").Value = 980
$ws.Range("H4").Value = 1.02
$ws.Range("I4").Value = 980
$ws.Range("J4").Value = 1.02
$ws.Range("K4").Value = 980
$ws.Range("N4").Value = 1.1
$ws.Range("P4").Value = 1.09
$ws.Range("Q4").Value = 1.01
$ws.Range("S4").Value = 1.26
$ws.Range("T4").Value = 1.04
$ws.Range("U4").Value = 1.04
$ws.Range("V4").Value = 1.01
$ws.Range("W4").Value = 1.01
$ws.Range("F5").Value = 5.4
$ws.Range("G5").Value = 5.5
$ws.Range("J5").Value = 3.7
$ws.Range("K5").Value = 3.75
$ws.Range("L5").Value = 1.45
$ws.Range("N5").Value = 3.35
$ws.Range("U5").Value = 1.9
$ws.Range("W5").Value = 1.22
$ws.Range("Z5").Value = 9.800000000000001
$ws.Range("AA5").Value = 18.5
$ws.Range("AM5").Value = 180
$ws.Range("AO5").Value = 15
$ws.Range("F6").Value = 1.02
$ws.Range("G6").Value = 980
$ws.Range("H6").Value = 1.02
$ws.Range("I6").Value = 980
$ws.Range("J6").Value = 1.02
$ws.Range("K6").Value = 980
$ws.Range("N6").Value = 1.16
$ws.Range("P6").Value = 1.16
$ws.Range("Q6").Value = 1.46
$ws.Range("R6").Value = 1.09
$ws.Range("S6").Value = 2.12
$ws.Range("T6").Value = 1.04
$ws.Range("U6").Value = 1.04
$ws.Range("V6").Value = 1.01
$ws.Range("W6").Value = 1.01
$ws.Range("F7").Value = 5.9
$ws.Range("G7").Value = 6
$ws.Range("H7").Value = 1.81
$ws.Range("I7").Value = 1.82
$ws.Range("L7").Value = 1.51
$ws.Range("P7").Value = 1.66
$ws.Range("T7").Value = 2.24
$ws.Range("V7").Value = 2.2
$ws.Range("W7").Value = 1.2
$ws.Range("AA7").Value = 18.5
$ws.Range("AF7").Value = 42
$ws.Range("AH7").Value = 27
$ws.Range("AO7").Value = 17
$ws.Range("F8").Value = 5.5
$ws.Range("G8").Value = 5.6
$ws.Range("I8").Value = 1.77
$ws.Range("J8").Value = 3.95
$ws.Range("K8").Value = 4
$ws.Range("L8").Value = 1.39
$ws.Range("W8").Value = 1.21
$ws.Range("Y8").Value = 7.8
$ws.Range("AH8").Value = 22
$ws.Range("AI8").Value = 38
$ws.Range("AJ8").Value = 150
$ws.Range("AL8").Value = 85
$ws.Range("F9").Value = 2.86
$ws.Range("G9").Value = 2.94
$ws.Range("H9").Value = 2.66
$ws.Range("I9").Value = 2.7
$ws.Range("J9").Value = 3.55
$ws.Range("K9").Value = 3.6
$ws.Range("N9").Value = 3.5
$ws.Range("O9").Value = 1.34
$ws.Range("Q9").Value = 2.04
$ws.Range("S9").Value = 3.85
$ws.Range("T9").Value = 1.83
$ws.Range("V9").Value = 1.58
$ws.Range("W9").Value = 1.52
$ws.Range("Y9").Value = 1000
$ws.Range("Z9").Value = 1000
$ws.Range("AA9").Value = 980
$ws.Range("AB9").Value = 1000
$ws.Range("AF9").Value = 1000
$ws.Range("AL9").Value = 980
$ws.Range("AM9").Value = 110
$ws.Range("AN9").Value = 32
$ws.Range("AO9").Value = 28
$ws.Range("P10").Value = 1.96
$ws.Range("Q10").Value = 2.02
$ws.Range("AC10").Value = 7.6
$ws.Range("AL10").Value = 50
$ws.Range("AM10").Value = 85
$ws.Range("F11").Value = 1.62
$ws.Range("G11").Value = 1.64
$ws.Range("H11").Value = 6.4
$ws.Range("I11").Value = 7.4
$ws.Range("K11").Value = 4.5
$ws.Range("N11").Value = 4.9
$ws.Range("P11").Value = 2.24
$ws.Range("U11").Value = 2.08
$ws.Range("W11").Value = 2.56
$ws.Range("X11").Value = 22
$ws.Range("Z11").Value = 1000
$ws.Range("AB11").Value = 11
$ws.Range("AC11").Value = 9.800000000000001
$ws.Range("AE11").Value = 1000
$ws.Range("AH11").Value = 21
$ws.Range("AI11").Value = 1000
$ws.Range("AN11").Value = 7.8
$ws.Range("AO11").Value = 1000
